# Add season-record columns (Wins / Losses / Ties) to the team stats sheet.
# Mirrors commit: "Created functions to get season record" — new columns
# AD:AF are appended after the existing AC ("Unnamed: 28") column, with a
# header row (styled like the rest of row 1) and the same W/L/T totals
# (88 / 74 / 0) repeated down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell onto the three new header cells, then write the
# header text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-47). Every team in this
# sheet shares the same season record, so the same three numbers repeat.
$lastRow = 47

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 88
    $ws.Cells.Item($row, 31).Value = 74
    $ws.Cells.Item($row, 32).Value = 0
}
